$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("flight")

$ws.Range("B6").Value = "Hellos"
$ws.Range("B6").Select()
